$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.246.14"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.260.17"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'308.37"
$ws.Range("E5").Value = "  -5.57%  "
$ws.Range("D6").Value = "'98.88"
$ws.Range("E6").Value = "  -4.76%  "
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  -4.64%  "
$ws.Range("D10").Value = "'36.03"
$ws.Range("E10").Value = "  -5.94%  "
$ws.Range("D11").Value = "'0.0824"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "'7.37"
$ws.Range("E12").Value = "  -6.00%  "
$ws.Range("D13").Value = "'0.105"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").Value = "2.604.24"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "'0.844"
$ws.Range("E15").Value = "  -3.30%  "
$ws.Range("D16").Value = "2.255.99"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'13.93"
$ws.Range("E17").Value = "  -3.68%  "
$ws.Range("D18").Value = "44.118.11"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'13.02"
$ws.Range("E19").Value = "  -8.28%  "
$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "'6.37"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Value = "'65.52"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'239.78"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  -9.76%  "
$ws.Range("D25").Value = "'2.00"
$ws.Range("E25").Value = "  -8.65%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'10.22"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'37.10"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("E29").Value = "  -4.38%  "
$ws.Range("D30").Value = "'6.13"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("D31").Value = "'20.27"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "'157.51"
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").Value = "'0.0835"
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("D34").Value = "'3.44"
$ws.Range("E34").Value = "  +11.31%  "
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").Value = "'1.90"
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("D37").Value = "'0.120"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  -7.19%  "
$ws.Range("D39").Value = "'16.23"
$ws.Range("E39").Value = "  +7.88%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'3.93"
$ws.Range("E40").Value = "  -9.58%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "'3.44"
$ws.Range("E41").Value = "  -11.68%  "
$ws.Range("D42").Value = "'0.0309"
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.738.25"
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "'88.20"
$ws.Range("E45").Value = "  +5.53%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'16.03"
$ws.Range("E46").Value = "  +10.52%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.195"
$ws.Range("E47").Value = "  -4.45%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'5.18"
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'102.37"
$ws.Range("E49").Value = "  -3.25%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "'70.98"
$ws.Range("E50").Value = "  -6.78%  "
$ws.Range("E51").Value = "  -2.70%  "
